# (JMT) Updated 1s6, 1s7, 1s11 to meet acceleration specification
#
# Inserts a new "npc_speed_increase" parameter column (with value "[0..30]kph")
# into the bl_1s11 parameter table, between the existing "npc_speed_gt_dut_start"
# (col E) and "npc_dist_gt_dut_start" (col F) columns, shifting the remaining
# columns one place to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F, pushing the old F..J columns to G..K.
$ws.Columns("F:F").Insert()

# Give the new column the same (bestFit-ish) width as its neighbour, column E,
# rather than leaving it at the sheet default.
$ws.Columns("F:F").ColumnWidth = 21.5

# Populate the new column's header (row 2) and value (row 3) cells.
# Write the data row first, then the header, so the shared-string table
# gets the two new entries in that same order.
$ws.Range("F3").Value = "[0..30]kph"
$ws.Range("F2").Value = "npc_speed_increase"

# Match the post-edit selection left behind in the saved workbook.
$null = $ws.Range("F2").Select()
